$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to retain a literal text value (avoids Excel
    # auto-converting numeric-looking strings like "316.80" into numbers),
    # then restore the default "Normal" style so no stray formatting is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "28.514.14"
$ws.Cells.Item(2, 5).Value = "  -0.39%  "
$ws.Cells.Item(3, 4).Value = "1.820.11"
$ws.Cells.Item(3, 5).Value = "  -0.47%  "
$ws.Cells.Item(4, 5).Value = "  +0.18%  "
Set-TextValue $ws.Cells.Item(5, 4) "316.80"
$ws.Cells.Item(5, 5).Value = "  +0.22%  "
$ws.Cells.Item(6, 5).Value = "  +0.16%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.5168"
$ws.Cells.Item(7, 5).Value = "  -3.28%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.3878"
$ws.Cells.Item(8, 5).Value = "  -2.92%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.08433"
$ws.Cells.Item(9, 5).Value = "  +8.37%  "
Set-TextValue $ws.Cells.Item(10, 4) "41.83"
$ws.Cells.Item(10, 5).Value = "  -0.41%  "
$ws.Cells.Item(11, 5).Value = "  -1.03%  "
Set-TextValue $ws.Cells.Item(12, 4) "6.436"
$ws.Cells.Item(12, 5).Value = "  +1.70%  "
Set-TextValue $ws.Cells.Item(13, 4) "21.02"
$ws.Cells.Item(13, 5).Value = "  -0.80%  "
Set-TextValue $ws.Cells.Item(14, 4) "1.002"
$ws.Cells.Item(14, 5).Value = "  +0.15%  "
Set-TextValue $ws.Cells.Item(15, 4) "7.512"
$ws.Cells.Item(15, 5).Value = "  -1.14%  "
$ws.Cells.Item(16, 4).Value = "1.819.86"
$ws.Cells.Item(16, 5).Value = "  -0.29%  "
$ws.Cells.Item(17, 5).Value = "  +3.79%  "
Set-TextValue $ws.Cells.Item(18, 4) "92.89"
$ws.Cells.Item(18, 5).Value = "  -0.27%  "
Set-TextValue $ws.Cells.Item(19, 4) "0.06677"
$ws.Cells.Item(19, 5).Value = "  +1.34%  "
Set-TextValue $ws.Cells.Item(20, 4) "17.73"
$ws.Cells.Item(20, 5).Value = "  -0.49%  "
$ws.Cells.Item(21, 5).Value = "  +0.11%  "
Set-TextValue $ws.Cells.Item(22, 4) "6.081"
$ws.Cells.Item(22, 5).Value = "  -0.25%  "
$ws.Cells.Item(23, 4).Value = "28.554.50"
$ws.Cells.Item(23, 5).Value = "  -0.24%  "
$ws.Cells.Item(24, 5).Value = "  +1.42%  "
Set-TextValue $ws.Cells.Item(25, 4) "2.270"
$ws.Cells.Item(25, 5).Value = "  +1.65%  "
Set-TextValue $ws.Cells.Item(26, 4) "21.06"
$ws.Cells.Item(26, 5).Value = "  +1.03%  "
Set-TextValue $ws.Cells.Item(27, 4) "158.98"
$ws.Cells.Item(27, 5).Value = "  +1.46%  "
$ws.Cells.Item(28, 4).Value = "2.031.78"
$ws.Cells.Item(28, 5).Value = "  -0.24%  "
Set-TextValue $ws.Cells.Item(29, 4) "2.415"
$ws.Cells.Item(29, 5).Value = "  -0.05%  "
Set-TextValue $ws.Cells.Item(30, 4) "126.01"
$ws.Cells.Item(30, 5).Value = "  +0.28%  "
Set-TextValue $ws.Cells.Item(31, 4) "0.1086"
$ws.Cells.Item(31, 5).Value = "  -3.41%  "
$ws.Cells.Item(32, 5).Value = "  -5.31%  "
$ws.Cells.Item(33, 5).Value = "  -0.47%  "
Set-TextValue $ws.Cells.Item(34, 4) "0.07506"
$ws.Cells.Item(34, 5).Value = "  +1.85%  "
Set-TextValue $ws.Cells.Item(35, 4) "3.684"
$ws.Cells.Item(35, 5).Value = "  +0.87%  "
Set-TextValue $ws.Cells.Item(36, 4) "0.2232"
$ws.Cells.Item(36, 5).Value = "  -1.81%  "
$ws.Cells.Item(37, 5).Value = "  +0.30%  "
Set-TextValue $ws.Cells.Item(38, 4) "5.203"
$ws.Cells.Item(38, 5).Value = "  -0.32%  "
Set-TextValue $ws.Cells.Item(39, 4) "8.752"
$ws.Cells.Item(39, 5).Value = "  -1.90%  "
$ws.Cells.Item(40, 5).Value = "  +0.23%  "
Set-TextValue $ws.Cells.Item(41, 4) "11.25"
$ws.Cells.Item(41, 5).Value = "  -1.47%  "
Set-TextValue $ws.Cells.Item(42, 4) "1.194"
$ws.Cells.Item(42, 5).Value = "  -0.21%  "
Set-TextValue $ws.Cells.Item(43, 4) "1.402"
$ws.Cells.Item(43, 5).Value = "  +0.57%  "
Set-TextValue $ws.Cells.Item(44, 4) "13.57"
$ws.Cells.Item(44, 5).Value = "  +0.14%  "
Set-TextValue $ws.Cells.Item(45, 4) "3.779"
$ws.Cells.Item(45, 5).Value = "  +1.83%  "
Set-TextValue $ws.Cells.Item(46, 4) "0.5942"
$ws.Cells.Item(46, 5).Value = "  -0.23%  "
Set-TextValue $ws.Cells.Item(47, 4) "126.03"
$ws.Cells.Item(47, 5).Value = "  +0.15%  "
Set-TextValue $ws.Cells.Item(48, 4) "1.991"
$ws.Cells.Item(48, 5).Value = "  -0.52%  "
Set-TextValue $ws.Cells.Item(49, 4) "1.201"
$ws.Cells.Item(49, 5).Value = "  +0.48%  "
Set-TextValue $ws.Cells.Item(50, 4) "0.06974"
$ws.Cells.Item(50, 5).Value = "  +0.10%  "
Set-TextValue $ws.Cells.Item(51, 4) "74.40"
$ws.Cells.Item(51, 5).Value = "  -0.40%  "
